$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "RGossF"

# Fix a tiny rounding difference in I13 (Gaussian quadrature recompute)
$ws.Range("I13").Value = 0.9953620408711823

# Append new row 16 with an additional Gaussian-quadrature averaged intensity entry
$ws.Range("A16").Value = 14
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.003347938170875
$ws.Range("D16").Value = 0.9733256173515474
$ws.Range("E16").Value = 1.001125728766453
$ws.Range("F16").Value = 1.003347938170875
$ws.Range("G16").Value = 0.9836962612942944
$ws.Range("H16").Value = 1.012687467361678
$ws.Range("I16").Value = 0.9998331504704113
$ws.Range("J16").Value = 0.9733256173515474
$ws.Range("K16").Value = 0.9872256730590003
$ws.Range("L16").Value = 0.9952868056149375
$ws.Range("M16").Value = 0.9956693605692098
